$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# A leading apostrophe forces Excel to store these as text (preserving
# formats like '25.759.00', '1.000', leading zeros, percent strings, etc.)

# Row 2
$ws.Range("D2").Value = "'25.759.00"
$ws.Range("E2").Value = "'  -5.35%  "

# Row 3
$ws.Range("D3").Value = "'1.814.94"
$ws.Range("E3").Value = "'  -4.68%  "

# Row 4
$ws.Range("E4").Value = "'  +0.09%  "

# Row 5
$ws.Range("D5").Value = "'276.75"
$ws.Range("E5").Value = "'  -9.62%  "

# Row 6
$ws.Range("E6").Value = "'  +0.07%  "

# Row 7
$ws.Range("D7").Value = "'0.5001"
$ws.Range("E7").Value = "'  -6.95%  "

# Row 8
$ws.Range("D8").Value = "'0.3496"
$ws.Range("E8").Value = "'  -8.28%  "

# Row 9
$ws.Range("D9").Value = "'44.11"
$ws.Range("E9").Value = "'  -3.98%  "

# Row 10
$ws.Range("D10").Value = "'0.06625"
$ws.Range("E10").Value = "'  -9.02%  "

# Row 11
$ws.Range("D11").Value = "'20.05"
$ws.Range("E11").Value = "'  -9.70%  "

# Row 12
$ws.Range("D12").Value = "'0.8436"
$ws.Range("E12").Value = "'  -6.77%  "

# Row 13
$ws.Range("D13").Value = "'0.07810"
$ws.Range("E13").Value = "'  -4.73%  "

# Row 14
$ws.Range("D14").Value = "'1.800.83"
$ws.Range("E14").Value = "'  +61.64%  "

# Row 15
$ws.Range("D15").Value = "'5.046"
$ws.Range("E15").Value = "'  -5.48%  "

# Row 16
$ws.Range("D16").Value = "'87.54"
$ws.Range("E16").Value = "'  -8.81%  "

# Row 17
$ws.Range("E17").Value = "'  +0.16%  "

# Row 18
$ws.Range("D18").Value = "'13.89"
$ws.Range("E18").Value = "'  -6.46%  "

# Row 19
$ws.Range("B19").Value = "'Dai"
$ws.Range("C19").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "'  +0.06%  "

# Row 20
$ws.Range("B20").Value = "'ShibaInu"
$ws.Range("C20").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007978"
$ws.Range("E20").Value = "'  -7.73%  "

# Row 21
$ws.Range("D21").Value = "'25.819.02"
$ws.Range("E21").Value = "'  -5.22%  "

# Row 22
$ws.Range("D22").Value = "'4.732"
$ws.Range("E22").Value = "'  -6.13%  "

# Row 23
$ws.Range("D23").Value = "'10.01"
$ws.Range("E23").Value = "'  -7.09%  "

# Row 24
$ws.Range("D24").Value = "'6.072"
$ws.Range("E24").Value = "'  -6.67%  "

# Row 25
$ws.Range("D25").Value = "'140.89"
$ws.Range("E25").Value = "'  -5.66%  "

# Row 26
$ws.Range("D26").Value = "'2.103"
$ws.Range("E26").Value = "'  -8.29%  "

# Row 27
$ws.Range("D27").Value = "'1.662"
$ws.Range("E27").Value = "'  -5.18%  "

# Row 28
$ws.Range("D28").Value = "'16.84"
$ws.Range("E28").Value = "'  -8.24%  "

# Row 29
$ws.Range("D29").Value = "'108.49"
$ws.Range("E29").Value = "'  -6.91%  "

# Row 30
$ws.Range("D30").Value = "'4.323"
$ws.Range("E30").Value = "'  -10.17%  "

# Row 31
$ws.Range("D31").Value = "'4.204"
$ws.Range("E31").Value = "'  -11.14%  "

# Row 32
$ws.Range("D32").Value = "'0.08753"
$ws.Range("E32").Value = "'  -4.98%  "

# Row 33
$ws.Range("D33").Value = "'0.04834"
$ws.Range("E33").Value = "'  -4.82%  "

# Row 34
$ws.Range("D34").Value = "'0.7374"
$ws.Range("E34").Value = "'  -11.30%  "

# Row 35
$ws.Range("B35").Value = "'HuobiToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.878"
$ws.Range("E35").Value = "'  -4.11%  "

# Row 36
$ws.Range("B36").Value = "'ARBITRUM"
$ws.Range("C36").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "'1.128"
$ws.Range("E36").Value = "'  -7.27%  "

# Row 37
$ws.Range("D37").Value = "'1.000"
$ws.Range("E37").Value = "'  -0.33%  "

# Row 38
$ws.Range("D38").Value = "'3.058"
$ws.Range("E38").Value = "'  -8.44%  "

# Row 39
$ws.Range("D39").Value = "'2.467"
$ws.Range("E39").Value = "'  -7.96%  "

# Row 40
$ws.Range("D40").Value = "'0.5300"
$ws.Range("E40").Value = "'  -9.23%  "

# Row 41
$ws.Range("D41").Value = "'0.01865"
$ws.Range("E41").Value = "'  -6.96%  "

# Row 42
$ws.Range("D42").Value = "'0.9698"
$ws.Range("E42").Value = "'  -9.86%  "

# Row 43
$ws.Range("D43").Value = "'111.96"
$ws.Range("E43").Value = "'  -4.43%  "

# Row 44
$ws.Range("D44").Value = "'6.221"
$ws.Range("E44").Value = "'  -5.98%  "

# Row 45
$ws.Range("D45").Value = "'8.148"
$ws.Range("E45").Value = "'  -12.55%  "

# Row 46
$ws.Range("D46").Value = "'0.4675"
$ws.Range("E46").Value = "'  -6.83%  "

# Row 47
$ws.Range("D47").Value = "'1.001"
$ws.Range("E47").Value = "'  +0.06%  "

# Row 48
$ws.Range("E48").Value = "'  -8.78%  "

# Row 49
$ws.Range("D49").Value = "'9.203"
$ws.Range("E49").Value = "'  -8.72%  "

# Row 50
$ws.Range("D50").Value = "'35.60"
$ws.Range("E50").Value = "'  -7.25%  "

# Row 51
$ws.Range("B51").Value = "'Cronos"
$ws.Range("C51").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05874"
$ws.Range("E51").Value = "'  -4.70%  "
